$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1395.0769
$ws.Range("I19").Value = 999.5714
$ws.Range("K19").Value = 999.5714
$ws.Range("M19").Value = -824.5714
$ws.Range("H62").Value = 8435.096
$ws.Range("I62").Value = 5821.3105
$ws.Range("J62").Value = 14265.846
$ws.Range("K62").Value = 5821.3105
$ws.Range("L62").Value = 14265.846
$ws.Range("M62").Value = -5197.3105
$ws.Range("N62").Value = -15513.846
$ws.Range("H65").Value = 8435.096
$ws.Range("I65").Value = 5821.3105
$ws.Range("J65").Value = 14265.846
$ws.Range("K65").Value = 29106.5525
$ws.Range("L65").Value = 71329.23
$ws.Range("M65").Value = -25986.5525
$ws.Range("N65").Value = -77569.23
$ws.Range("H116").Value = 6922849.5
$ws.Range("J116").Value = 1998.75
$ws.Range("L116").Value = 1998.75
$ws.Range("N116").Value = -8882.75
$ws.Range("H137").Value = 27028368
$ws.Range("I137").Value = 43479296
$ws.Range("K137").Value = 130437888
$ws.Range("M137").Value = -130435338

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21464.436
$ws.Range("I32").Value = 3080
$ws.Range("K32").Value = 3080
$ws.Range("M32").Value = -2793
$ws.Range("H45").Value = 897.3333
$ws.Range("I45").Value = 875.4
$ws.Range("J45").Value = 1007
$ws.Range("K45").Value = 875.4
$ws.Range("L45").Value = 1007
$ws.Range("M45").Value = -498.4
$ws.Range("N45").Value = -1761
$ws.Range("H61").Value = 1645.6
$ws.Range("I61").Value = 1110.88
$ws.Range("J61").Value = 4319.2
$ws.Range("K61").Value = 1110.88
$ws.Range("L61").Value = 4319.2
$ws.Range("M61").Value = -898.8800000000001
$ws.Range("N61").Value = -4743.2
$ws.Range("H74").Value = 4188.9023
$ws.Range("I74").Value = 1012.7143
$ws.Range("J74").Value = 22716.666
$ws.Range("K74").Value = 1012.7143
$ws.Range("L74").Value = 22716.666
$ws.Range("M74").Value = -138.7143
$ws.Range("N74").Value = -24464.666
$ws.Range("H77").Value = 4188.9023
$ws.Range("I77").Value = 1012.7143
$ws.Range("J77").Value = 22716.666
$ws.Range("K77").Value = 5063.5715
$ws.Range("L77").Value = 113583.33
$ws.Range("M77").Value = -695.5715
$ws.Range("N77").Value = -122319.33
$ws.Range("H88").Value = 4725.75
$ws.Range("I88").Value = 2449.5
$ws.Range("J88").Value = 5484.5
$ws.Range("K88").Value = 2449.5
$ws.Range("L88").Value = 5484.5
$ws.Range("M88").Value = -2043.5
$ws.Range("N88").Value = -6296.5
$ws.Range("H91").Value = 4725.75
$ws.Range("I91").Value = 2449.5
$ws.Range("J91").Value = 5484.5
$ws.Range("K91").Value = 2449.5
$ws.Range("L91").Value = 5484.5
$ws.Range("M91").Value = -1045.5
$ws.Range("N91").Value = -8292.5
$ws.Range("H132").Value = 2458.1226
$ws.Range("I132").Value = 2049.4119
$ws.Range("J132").Value = 3384.5334
$ws.Range("K132").Value = 6148.2357
$ws.Range("L132").Value = 10153.6002
$ws.Range("M132").Value = -3618.2357
$ws.Range("N132").Value = -15213.6002
$ws.Range("H136").Value = 1645.6
$ws.Range("I136").Value = 1110.88
$ws.Range("J136").Value = 4319.2
$ws.Range("K136").Value = 3332.64
$ws.Range("L136").Value = 12957.6
$ws.Range("M136").Value = -782.6400000000003
$ws.Range("N136").Value = -18057.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 12900.1
$ws.Range("I86").Value = 3798.8
$ws.Range("J86").Value = 22001.4
$ws.Range("K86").Value = 3798.8
$ws.Range("L86").Value = 22001.4
$ws.Range("M86").Value = -2675.8
$ws.Range("N86").Value = -24247.4
$ws.Range("H89").Value = 12900.1
$ws.Range("I89").Value = 3798.8
$ws.Range("J89").Value = 22001.4
$ws.Range("K89").Value = 18994
$ws.Range("L89").Value = 110007
$ws.Range("M89").Value = -13378
$ws.Range("N89").Value = -121239

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1977.3334
$ws.Range("I31").Value = 1060.421
$ws.Range("J31").Value = 3221.7144
$ws.Range("K31").Value = 1060.421
$ws.Range("L31").Value = 3221.7144
$ws.Range("M31").Value = -765.421
$ws.Range("N31").Value = -3811.7144
$ws.Range("H34").Value = 1977.3334
$ws.Range("I34").Value = 1060.421
$ws.Range("J34").Value = 3221.7144
$ws.Range("K34").Value = 1060.421
$ws.Range("L34").Value = 3221.7144
$ws.Range("M34").Value = -858.421
$ws.Range("N34").Value = -3625.7144
$ws.Range("H62").Value = 21763.545
$ws.Range("I62").Value = 31742.715
$ws.Range("K62").Value = 31742.715
$ws.Range("M62").Value = -31118.715
$ws.Range("H65").Value = 21763.545
$ws.Range("I65").Value = 31742.715
$ws.Range("K65").Value = 158713.575
$ws.Range("M65").Value = -155593.575
$ws.Range("H133").Value = 26949.938
$ws.Range("J133").Value = 29013.04
$ws.Range("L133").Value = 29013.04
$ws.Range("N133").Value = -34073.04
$ws.Range("H135").Value = 39705.355
$ws.Range("J135").Value = 39705.355
$ws.Range("L135").Value = 39705.355
$ws.Range("N135").Value = -49845.355

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 624.9167
$ws.Range("I46").Value = 99.666664
$ws.Range("K46").Value = 298.999992
$ws.Range("M46").Value = -207.999992

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6426.448
$ws.Range("I70").Value = 6727.6113
$ws.Range("K70").Value = 6727.6113
$ws.Range("M70").Value = -6457.6113
$ws.Range("H73").Value = 6426.448
$ws.Range("I73").Value = 6727.6113
$ws.Range("K73").Value = 6727.6113
$ws.Range("M73").Value = -5791.6113

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 50000
$ws.Range("I33").Value = 50000
$ws.Range("K33").Value = 50000
$ws.Range("M33").Value = -49710
$ws.Range("H43").Value = 18250
$ws.Range("J43").Value = 7666.6665
$ws.Range("L43").Value = 7666.6665
$ws.Range("N43").Value = -8052.6665
$ws.Range("H55").Value = 308.36365
$ws.Range("I55").Value = 260.2857
$ws.Range("J55").Value = 330.8
$ws.Range("K55").Value = 260.2857
$ws.Range("L55").Value = 330.8
$ws.Range("M55").Value = -87.28570000000002
$ws.Range("N55").Value = -676.8
$ws.Range("H104").Value = 30400.2
$ws.Range("J104").Value = 30400.2
$ws.Range("L104").Value = 30400.2
$ws.Range("N104").Value = -37388.2
$ws.Range("H139").Value = 48944.445
$ws.Range("J139").Value = 48944.445
$ws.Range("L139").Value = 48944.445
$ws.Range("N139").Value = -59224.445

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 60000
$ws.Range("J139").Value = 60000
$ws.Range("L139").Value = 60000
$ws.Range("N139").Value = -70280
